$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 72, shifting existing rows 72:164 down to 73:165
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new data point
$ws.Cells.Item(72, 1).Value = 10
$ws.Cells.Item(72, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(72, 3).Value = "La Araucanía"
$ws.Cells.Item(72, 4).Value = 44413
$ws.Cells.Item(72, 5).Value = 9
$ws.Cells.Item(72, 6).Value = 100114013
$ws.Cells.Item(72, 7).Value = "Zanahoria"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 200
$ws.Cells.Item(72, 11).Value = 5000
$ws.Cells.Item(72, 12).Value = 5000
$ws.Cells.Item(72, 13).Value = 5000
$ws.Cells.Item(72, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(72, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(72, 16).Value = 200
$ws.Cells.Item(72, 17).Value = 25
$ws.Cells.Item(72, 18).Value = "Hortaliza"

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
